$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the contents of columns B..AC (2..29) between two rows, while
# leaving column A (the running index) untouched on both rows.
# ---------------------------------------------------------------------------
function Swap-Rows($row1, $row2) {
    for ($col = 2; $col -le 29; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Matches in the re-scraped source swapped position with their neighbour row
Swap-Rows 9 10
Swap-Rows 29 30
Swap-Rows 87 88

# ---------------------------------------------------------------------------
# Append the newly scraped match as row 146
# ---------------------------------------------------------------------------
$newRow = 146

# Copy formatting (style) only from the row above for the styled columns,
# so no new cellXf entries are created.
$ws.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow - 1, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 144
$ws.Cells.Item($newRow, 2).Value = 7952743
$ws.Cells.Item($newRow, 3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item($newRow, 4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item($newRow, 5).Value = 45388.65625
$ws.Cells.Item($newRow, 6).Value = "FK Sarajevo"
$ws.Cells.Item($newRow, 7).Value = "NK Igman Konjic"

$ws.Cells.Item($newRow, 11).Value = 1.222
$ws.Cells.Item($newRow, 12).Value = 5.2
$ws.Cells.Item($newRow, 13).Value = 9.5
$ws.Cells.Item($newRow, 14).Value = 1.181
$ws.Cells.Item($newRow, 15).Value = 5.5
$ws.Cells.Item($newRow, 16).Value = 12
$ws.Cells.Item($newRow, 17).Value = -2
$ws.Cells.Item($newRow, 18).Value = 2.025
$ws.Cells.Item($newRow, 19).Value = 1.775
$ws.Cells.Item($newRow, 20).Value = 3
$ws.Cells.Item($newRow, 21).Value = 1.875
$ws.Cells.Item($newRow, 22).Value = 1.925
$ws.Cells.Item($newRow, 23).Value = 0
$ws.Cells.Item($newRow, 24).Value = 0
$ws.Cells.Item($newRow, 25).Value = 0
$ws.Cells.Item($newRow, 26).Value = 0
$ws.Cells.Item($newRow, 27).Value = 0
